# Update crypto price/volume figures on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "27.414.13"; E = "  -0.55%  " },
    @{ Row = 3;  D = "1.562.45";  E = "  -1.05%  " },
    @{ Row = 4;  D = $null;       E = "  -0.05%  " },
    @{ Row = 5;  D = "208.02";    E = "  +0.99%  " },
    @{ Row = 6;  D = $null;       E = "  -0.62%  " },
    @{ Row = 7;  D = $null;       E = "  -0.04%  " },
    @{ Row = 8;  D = "21.79";     E = "  -1.33%  " },
    @{ Row = 9;  D = $null;       E = "  -1.23%  " },
    @{ Row = 10; D = $null;       E = "  +0.10%  " },
    @{ Row = 11; D = "0.0867";    E = "  +0.16%  " },
    @{ Row = 12; D = "1.785.47";  E = "  -0.92%  " },
    @{ Row = 13; D = "1.567.54";  E = "  -0.84%  " },
    @{ Row = 14; D = "3.81";      E = "  -0.72%  " },
    @{ Row = 15; D = "0.514";     E = "  -2.16%  " },
    @{ Row = 16; D = "63.25";     E = "  +0.88%  " },
    @{ Row = 17; D = "27.419.91"; E = "  -0.41%  " },
    @{ Row = 18; D = "212.47";    E = "  -2.37%  " },
    @{ Row = 19; D = $null;       E = "  -0.48%  " },
    @{ Row = 20; D = "7.23";      E = "  -1.00%  " },
    @{ Row = 21; D = $null;       E = "  -0.07%  " },
    @{ Row = 22; D = "4.10";      E = "  -0.45%  " },
    @{ Row = 23; D = "9.51";      E = "  +1.02%  " },
    @{ Row = 24; D = "2.01";      E = "  +1.37%  " },
    @{ Row = 25; D = "152.80";    E = "  -0.26%  " },
    @{ Row = 26; D = $null;       E = "  -0.01%  " },
    @{ Row = 27; D = $null;       E = "  +0.87%  " },
    @{ Row = 28; D = "14.97";     E = "  -0.23%  " },
    @{ Row = 29; D = $null;       E = "  -1.56%  " },
    @{ Row = 30; D = $null;       E = "  +0.03%  " },
    @{ Row = 31; D = $null;       E = "  +1.77%  " },
    @{ Row = 32; D = "3.18";      E = "  -1.36%  " },
    @{ Row = 33; D = "1.359.36";  E = "  -0.29%  " },
    @{ Row = 34; D = $null;       E = "  +0.08%  " },
    @{ Row = 35; D = "1.52";      E = "  +1.21%  " },
    @{ Row = 36; D = "0.974";     E = "  +1.30%  " },
    @{ Row = 37; D = $null;       E = "  +0.07%  " },
    @{ Row = 38; D = $null;       E = "  +1.89%  " },
    @{ Row = 39; D = "0.531";     E = "  -0.33%  " },
    @{ Row = 40; D = "0.819";     E = "  +0.96%  " },
    @{ Row = 41; D = $null;       E = "  -0.03%  " },
    @{ Row = 42; D = $null;       E = "  -0.04%  " },
    @{ Row = 43; D = $null;       E = "  +1.87%  " },
    @{ Row = 44; D = "64.01";     E = "  +1.17%  " },
    @{ Row = 45; D = $null;       E = "  +0.52%  " },
    @{ Row = 47; D = "1.698.09";  E = "  -0.72%  " },
    @{ Row = 48; D = "85.40";     E = "  -2.08%  " },
    @{ Row = 49; D = "0.0₇0984"; E = "  -1.57%  " },
    @{ Row = 50; D = $null;       E = "  -1.37%  " },
    @{ Row = 51; D = $null;       E = "  -0.49%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cellD = $ws.Cells.Item($u.Row, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.D
    }
    $cellE = $ws.Cells.Item($u.Row, 5)
    $cellE.NumberFormat = "@"
    $cellE.Value = $u.E
}
